$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 894, pushing existing rows 894-935 down to 895-936
$ws.Rows(894).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real Excel
# date value. Force text number-format before assigning so Excel doesn't
# auto-convert the "yyyy/mm/dd"-looking string into a date serial, then
# restore the default "Normal" style so no stray explicit style id lingers
# on the cell (matching the rest of the data rows, which carry no style).
$ws.Cells.Item(894, 1).NumberFormat = "@"
$ws.Cells.Item(894, 1).Value = "2026/02/28"
$ws.Cells.Item(894, 1).Style = "Normal"

$ws.Cells.Item(894, 2).Value = "土"
$ws.Cells.Item(894, 3).Value = 19
$ws.Cells.Item(894, 4).Value = 201
